$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data rows (Result, Date, Execute, PaymentType, TaxType)
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Fri Sep 08 18:12:13 EDT 2023"
$ws.Range("C2").Value = "Y"
$ws.Range("D2").Value = "Quarterly Estimated Tax"
$ws.Range("E2").Value = "Personal Income Tax"

$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Fri Sep 08 18:12:29 EDT 2023"
$ws.Range("C3").Value = "Y"
$ws.Range("D3").Value = "Extension Payments"
$ws.Range("E3").Value = "Personal Income Tax"

$ws.Range("A4").Value = "Pass"
$ws.Range("B4").Value = "Fri Sep 08 18:12:44 EDT 2023"
$ws.Range("C4").Value = "Y"
$ws.Range("D4").Value = "New Tax Return Amount Due"
$ws.Range("E4").Value = "Personal Income Tax"

# Update the active selection to D3
$ws.Range("D3").Select()
